$d = $word.ActiveDocument

# The technical-requirements "note" paragraph is made up of 5 runs, all
# bold. We need to: (a) lowercase the leading "Note:" to "note:", and
# (b) swap bold for italic on every run in the paragraph, while leaving
# everything else (run/paragraph boundaries, the rest of the document)
# untouched.
foreach ($para in $d.Paragraphs) {
    $full = $para.Range
    if ($full.Text -like "*Note: all tech requirements are subject to available equipment*") {
        # Exclude the trailing paragraph-mark character so the pPr (and the
        # paragraph itself) survive the replace; only the run content changes.
        $target = $d.Range($full.Start, $full.End - 1)

        $xml = @'
<?xml version="1.0" standalone="yes"?>
<?mso-application progid="Word.Document"?>
<pkg:package xmlns:pkg="http://schemas.microsoft.com/office/2006/xmlPackage"><pkg:part pkg:name="/word/document.xml" pkg:contentType="application/vnd.openxmlformats-officedocument.wordprocessingml.document.main+xml"><pkg:xmlData><w:document xmlns:w="http://schemas.openxmlformats.org/wordprocessingml/2006/main"><w:body><w:p><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">note: all tech requirements are subject to available equipment and</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">venues&#8212;workshops will (in general) take place in typical university classrooms</w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve"> </w:t></w:r><w:r><w:rPr><w:i/></w:rPr><w:t xml:space="preserve">with projector and stereo PA system</w:t></w:r></w:p></w:body></w:document></pkg:xmlData></pkg:part></pkg:package>
'@

        $target.InsertXML($xml)
    }
}
